$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the forecast column (C) values for rows 2-21 (dates 2025/11/25 down to 2025/10/29)
# These are refreshed forecast numbers from the latest model run.
$ws.Cells.Item(2,3).Value  = 828.271
$ws.Cells.Item(3,3).Value  = 824.631
$ws.Cells.Item(4,3).Value  = 823.0599999999999
$ws.Cells.Item(5,3).Value  = 813.8915
$ws.Cells.Item(6,3).Value  = 815.1398
$ws.Cells.Item(7,3).Value  = 815.9267
$ws.Cells.Item(8,3).Value  = 834.2756000000001
$ws.Cells.Item(9,3).Value  = 822.5669
$ws.Cells.Item(10,3).Value = 807.8872
$ws.Cells.Item(11,3).Value = 829.0517
$ws.Cells.Item(12,3).Value = 816.1226
$ws.Cells.Item(13,3).Value = 815.2746
$ws.Cells.Item(14,3).Value = 809.5371
$ws.Cells.Item(15,3).Value = 804.4885
$ws.Cells.Item(16,3).Value = 801.9944
$ws.Cells.Item(17,3).Value = 799.0966
$ws.Cells.Item(18,3).Value = 795.7944
$ws.Cells.Item(19,3).Value = 793.9400000000001
$ws.Cells.Item(20,3).Value = 786.8131
$ws.Cells.Item(21,3).Value = 792.5

# Row 22 shifts to the next trading day (2025/10/29) with its refreshed forecast value
# (dates are stored as plain text in this sheet, so force a text format to stop
# Excel from auto-converting the string into a date serial number, then restore
# the default "Normal" style so no stray number format sticks around)
$ws.Cells.Item(22,1).NumberFormat = "@"
$ws.Cells.Item(22,1).Value = "2025/10/29"
$ws.Cells.Item(22,1).Style = "Normal"
$ws.Cells.Item(22,3).Value = 786.2939

# Row 23: a new day of actual data arrived (2025/10/28) - now holds an actual
# value in column B instead of a forecast in column C
$ws.Cells.Item(23,1).NumberFormat = "@"
$ws.Cells.Item(23,1).Value = "2025/10/28"
$ws.Cells.Item(23,1).Style = "Normal"
$ws.Cells.Item(23,2).Value = 792.5
$ws.Cells.Item(23,3).ClearContents()

# Row 24: likewise shifts to 2025/10/27 and becomes an actual value in column B
$ws.Cells.Item(24,1).NumberFormat = "@"
$ws.Cells.Item(24,1).Value = "2025/10/27"
$ws.Cells.Item(24,1).Style = "Normal"
$ws.Cells.Item(24,2).Value = 786.5
$ws.Cells.Item(24,3).ClearContents()
